$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, shifting existing rows 21+ down by one
$ws.Rows.Item(21).Insert()

# Fill in the new row 21 with data
$ws.Range("A21").Value = 3
$ws.Range("B21").Value = "Femacal de La Calera"
$ws.Range("C21").Value = "Coquimbo"
$ws.Range("D21").Value = 44546
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 100112022
$ws.Range("G21").Value = "Arveja Verde"
$ws.Range("H21").Value = "Perfection"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 75
$ws.Range("K21").Value = 18000
$ws.Range("L21").Value = 18500
$ws.Range("M21").Value = 18267
$ws.Range("N21").Value = "`$/saco 25 kilos"
$ws.Range("O21").Value = "Provincia de Limarí"
$ws.Range("P21").Value = 731
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"
